# Updated symbol list on Wed Dec 21 15:48:39 UTC 2022 with GitHub Actions
# Refresh the crypto price/volume snapshot on Sheet1.
# Numeric-looking "Price" column values are stored as literal text in the
# source data (inline strings), so a leading apostrophe is used to force
# Excel to keep them as text instead of auto-converting to Double (which
# would introduce floating point noise like 248.69999999999999).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.70"
$ws.Range("D3").Value = "'22.58"
$ws.Range("D4").Value = "'5.221"
$ws.Range("D5").Value = "'0.05694"
$ws.Range("D6").Value = "'3.409"
$ws.Range("D7").Value = "'6.334"
$ws.Range("D8").Value = "'0.8077"
$ws.Range("D9").Value = "'0.8872"
$ws.Range("D11").Value = "'0.07439"
$ws.Range("D12").Value = "'0.03056"
$ws.Range("D13").Value = "'0.03100"
$ws.Range("D14").Value = "'0.09398"
$ws.Range("D15").Value = "'3.874"
$ws.Range("D17").Value = "'0.04796"
$ws.Range("D18").Value = "'0.01827"
$ws.Range("D19").Value = "'0.0005796"
$ws.Range("E19").Value = "18OneONEWorstin24h"
$ws.Range("D20").Value = "'0.006427"
$ws.Range("D21").Value = "'0.004988"
$ws.Range("D22").Value = "'0.0009961"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("D24").Value = "'3.687"
$ws.Range("D25").Value = "'2.165"
$ws.Range("D26").Value = "'0.3270"
$ws.Range("D27").Value = "'0.1368"
$ws.Range("D40").Value = "'0.03979"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006811"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002806"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.007810"
$ws.Range("D45").Value = "'0.00005586"
$ws.Range("D47").Value = "'0.4986"
$ws.Range("D48").Value = "'0.2050"
$ws.Range("D49").Value = "'0.00002100"
